$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 33341326
$ws.Range("I74").Value = 50006000
$ws.Range("K74").Value = 50006000
$ws.Range("M74").Value = -50005064
$ws.Range("H77").Value = 33341326
$ws.Range("I77").Value = 50006000
$ws.Range("K77").Value = 250030000
$ws.Range("M77").Value = -250025320
$ws.Range("H106").Value = 2261.3333
$ws.Range("I106").Value = 2261.3333
$ws.Range("K106").Value = 2261.3333
$ws.Range("M106").Value = -1630.3333
$ws.Range("H125").Value = 29415156
$ws.Range("I125").Value = 45456516
$ws.Range("J125").Value = 5998.1665
$ws.Range("K125").Value = 409108644
$ws.Range("L125").Value = 53983.4985
$ws.Range("M125").Value = -409106184
$ws.Range("N125").Value = -58903.4985
$ws.Range("H138").Value = 2861189
$ws.Range("I138").Value = 2626.1538
$ws.Range("K138").Value = 7878.4614
$ws.Range("M138").Value = -2738.4614

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4038571.8
$ws.Range("I32").Value = 4469347
$ws.Range("J32").Value = 18000.334
$ws.Range("K32").Value = 4469347
$ws.Range("L32").Value = 18000.334
$ws.Range("M32").Value = -4469060
$ws.Range("N32").Value = -18574.334
$ws.Range("H45").Value = 4482.55
$ws.Range("I45").Value = 2025.6666
$ws.Range("K45").Value = 2025.6666
$ws.Range("M45").Value = -1648.6666
$ws.Range("H61").Value = 9266.24
$ws.Range("I61").Value = 3591.4443
$ws.Range("J61").Value = 12458.3125
$ws.Range("K61").Value = 3591.4443
$ws.Range("L61").Value = 12458.3125
$ws.Range("M61").Value = -3379.4443
$ws.Range("N61").Value = -12882.3125
$ws.Range("H122").Value = 2508.4482
$ws.Range("I122").Value = 1928.2941
$ws.Range("K122").Value = 5784.8823
$ws.Range("M122").Value = -3334.8823
$ws.Range("H132").Value = 1168521
$ws.Range("I132").Value = 2004643.5
$ws.Range("J132").Value = 7239.6113
$ws.Range("K132").Value = 6013930.5
$ws.Range("L132").Value = 21718.8339
$ws.Range("M132").Value = -6011400.5
$ws.Range("N132").Value = -26778.8339
$ws.Range("H136").Value = 9266.24
$ws.Range("I136").Value = 3591.4443
$ws.Range("J136").Value = 12458.3125
$ws.Range("K136").Value = 10774.3329
$ws.Range("L136").Value = 37374.9375
$ws.Range("M136").Value = -8224.332900000001
$ws.Range("N136").Value = -42474.9375

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3463.5
$ws.Range("I105").Value = 2332.25
$ws.Range("K105").Value = 2332.25
$ws.Range("M105").Value = -585.25
$ws.Range("H134").Value = 9908.546
$ws.Range("I134").Value = 2999.1667
$ws.Range("J134").Value = 18199.8
$ws.Range("K134").Value = 8997.500100000001
$ws.Range("L134").Value = 54599.39999999999
$ws.Range("M134").Value = -6462.500100000001
$ws.Range("N134").Value = -59669.39999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 100034710
$ws.Range("H31").Value = 7919.394
$ws.Range("I31").Value = 3021.1765
$ws.Range("K31").Value = 3021.1765
$ws.Range("M31").Value = -2726.1765
$ws.Range("H33").Value = 4782.5
$ws.Range("I33").Value = 2000
$ws.Range("J33").Value = 7565
$ws.Range("K33").Value = 2000
$ws.Range("L33").Value = 7565
$ws.Range("M33").Value = -1621
$ws.Range("N33").Value = -8323
$ws.Range("H34").Value = 7919.394
$ws.Range("I34").Value = 3021.1765
$ws.Range("K34").Value = 3021.1765
$ws.Range("M34").Value = -2819.1765
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").Value = $null
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").Value = $null
$ws.Range("H107").Value = 1162.2858
$ws.Range("I107").Value = 519
$ws.Range("J107").Value = 3220.8
$ws.Range("K107").Value = 519
$ws.Range("L107").Value = 3220.8
$ws.Range("M107").Value = 1401
$ws.Range("N107").Value = -7060.8
$ws.Range("H132").Value = 13427.857
$ws.Range("I132").Value = 8000
$ws.Range("J132").Value = 14332.5
$ws.Range("K132").Value = 24000
$ws.Range("L132").Value = 42997.5
$ws.Range("M132").Value = -21470
$ws.Range("N132").Value = -48057.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 191.2
$ws.Range("I23").Value = 144.66667
$ws.Range("K23").Value = 434.00001
$ws.Range("M23").Value = -199.00001
$ws.Range("H68").Value = 2009.1034
$ws.Range("I68").Value = 1212.9
$ws.Range("J68").Value = 2428.158
$ws.Range("K68").Value = 3638.7
$ws.Range("L68").Value = 7284.474
$ws.Range("M68").Value = -2827.7
$ws.Range("N68").Value = -8906.474
$ws.Range("H71").Value = 2009.1034
$ws.Range("I71").Value = 1212.9
$ws.Range("J71").Value = 2428.158
$ws.Range("K71").Value = 10916.1
$ws.Range("L71").Value = 21853.422
$ws.Range("M71").Value = -6860.1
$ws.Range("N71").Value = -29965.422
$ws.Range("H113").Value = 1755.8
$ws.Range("J113").Value = 2289.3635
$ws.Range("L113").Value = 6868.0905
$ws.Range("N113").Value = -11208.0905
$ws.Range("H122").Value = 726363.8
$ws.Range("I122").Value = 1768787.5
$ws.Range("J122").Value = 1199.5217
$ws.Range("K122").Value = 15919087.5
$ws.Range("L122").Value = 10795.6953
$ws.Range("M122").Value = -15916637.5
$ws.Range("N122").Value = -15695.6953
$ws.Range("H137").Value = 112177
$ws.Range("J137").Value = 101668.91
$ws.Range("L137").Value = 305006.73
$ws.Range("N137").Value = -315206.73

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 5051482.5
$ws.Range("I10").Value = 7575001.5
$ws.Range("J10").Value = 4444
$ws.Range("K10").Value = 7575001.5
$ws.Range("L10").Value = 4444
$ws.Range("M10").Value = -7574832.5
$ws.Range("N10").Value = -4782
$ws.Range("H12").Value = 500
$ws.Range("I12").Value = 500
$ws.Range("K12").Value = 500
$ws.Range("M12").Value = -360
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = $null
$ws.Range("H52").Value = 89999.2
$ws.Range("J52").Value = 89999.2
$ws.Range("L52").Value = 89999.2
$ws.Range("N52").Value = -90517.2
$ws.Range("H132").Value = 4287.0386
$ws.Range("I132").Value = 2517.375
$ws.Range("K132").Value = 7552.125
$ws.Range("M132").Value = -5022.125

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3315.4443
$ws.Range("I16").Value = 3229.9375
$ws.Range("K16").Value = 3229.9375
$ws.Range("M16").Value = -3059.9375
$ws.Range("H20").Value = 1160000
$ws.Range("J20").Value = 1250000
$ws.Range("L20").Value = 1250000
$ws.Range("N20").Value = -1250452
$ws.Range("H22").Value = 1726.5555
$ws.Range("I22").Value = 869.46155
$ws.Range("K22").Value = 869.46155
$ws.Range("M22").Value = -574.46155
$ws.Range("H23").Value = 5966.6665
$ws.Range("I23").Value = 5000
$ws.Range("J23").Value = 6450
$ws.Range("K23").Value = 5000
$ws.Range("L23").Value = 6450
$ws.Range("M23").Value = -4770
$ws.Range("N23").Value = -6910
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").Value = $null
$ws.Range("H27").Value = 1726.5555
$ws.Range("I27").Value = 869.46155
$ws.Range("K27").Value = 869.46155
$ws.Range("M27").Value = -762.46155
$ws.Range("H40").Value = 55560030
$ws.Range("I40").Value = 55560030
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 55560030
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -55559894
$ws.Range("N40").Value = $null
$ws.Range("H55").Value = 917.5484
$ws.Range("I55").Value = 962.3158
$ws.Range("J55").Value = 846.6667
$ws.Range("K55").Value = 962.3158
$ws.Range("L55").Value = 846.6667
$ws.Range("M55").Value = -789.3158
$ws.Range("N55").Value = -1192.6667
$ws.Range("H61").Value = 8339124.5
$ws.Range("I61").Value = 25002500
$ws.Range("K61").Value = 25002500
$ws.Range("M61").Value = -25002298
$ws.Range("H104").Value = 39846
$ws.Range("J104").Value = 39846
$ws.Range("L104").Value = 39846
$ws.Range("N104").Value = -46834
$ws.Range("H113").Value = 8339124.5
$ws.Range("I113").Value = 25002500
$ws.Range("K113").Value = 25002500
$ws.Range("M113").Value = -25000330
$ws.Range("H132").Value = 6499.5405
$ws.Range("I132").Value = 3770.4707
$ws.Range("J132").Value = 8819.25
$ws.Range("K132").Value = 11311.4121
$ws.Range("L132").Value = 26457.75
$ws.Range("M132").Value = -8781.4121
$ws.Range("N132").Value = -31517.75
$ws.Range("H136").Value = 8020.0186
$ws.Range("J136").Value = 12253.05
$ws.Range("L136").Value = 36759.14999999999
$ws.Range("N136").Value = -41859.14999999999
$ws.Range("H141").Value = 74686.125
$ws.Range("J141").Value = 74686.125
$ws.Range("L141").Value = 74686.125
$ws.Range("N141").Value = -85046.125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 299999
$ws.Range("J41").Value = 299999
$ws.Range("L41").Value = 299999
$ws.Range("N41").Value = -300779
$ws.Range("H81").Value = 15005470
$ws.Range("I81").Value = 834715.4399999999
$ws.Range("K81").Value = 1669430.88
$ws.Range("M81").Value = -1668369.88
$ws.Range("H84").Value = 15005470
$ws.Range("I84").Value = 834715.4399999999
$ws.Range("K84").Value = 8347154.399999999
$ws.Range("M84").Value = -8341850.399999999
$ws.Range("H100").Value = 1480.9
$ws.Range("I100").Value = 545
$ws.Range("K100").Value = 1090
$ws.Range("M100").Value = -549
$ws.Range("H136").Value = 33609.293
$ws.Range("I136").Value = 1533.909
$ws.Range("K136").Value = 4601.727000000001
$ws.Range("M136").Value = -2051.727000000001
